$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Durations_psplib")

$ws.Range("C1").Value = 2441.530424594879

$ws.Range("A4").Value = 304051.1835
$ws.Range("B4").Value = 299190
$ws.Range("F4").Value = 102243.043
$ws.Range("G4").Value = 100450
